$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.68
$ws.Cells.Item(2, 9).Value = 9.6
$ws.Cells.Item(2, 11).Value = 7.4
$ws.Cells.Item(2, 12).Value = 1.26
$ws.Cells.Item(2, 13).Value = 1.03
$ws.Cells.Item(2, 14).Value = 4.7
$ws.Cells.Item(2, 16).Value = 2.26
$ws.Cells.Item(2, 17).Value = 1.62
$ws.Cells.Item(2, 18).Value = 1.5
$ws.Cells.Item(2, 19).Value = 2.58
$ws.Cells.Item(2, 20).Value = 1.68
$ws.Cells.Item(2, 21).Value = 1.87
$ws.Cells.Item(2, 22).Value = 1.11
$ws.Cells.Item(2, 23).Value = 2.46
$ws.Cells.Item(2, 24).Value = 32
$ws.Cells.Item(2, 25).Value = 40
$ws.Cells.Item(2, 26).Value = 95
$ws.Cells.Item(2, 28).Value = 14.5
$ws.Cells.Item(2, 29).Value = 16
$ws.Cells.Item(2, 30).Value = 40
$ws.Cells.Item(2, 32).Value = 14.5
$ws.Cells.Item(2, 36).Value = 20
$ws.Cells.Item(2, 38).Value = 46
$ws.Cells.Item(3, 6).Value = 2.6
$ws.Cells.Item(3, 7).Value = 3.25
$ws.Cells.Item(3, 9).Value = 2.7
$ws.Cells.Item(3, 10).Value = 3.25
$ws.Cells.Item(3, 14).Value = 4
$ws.Cells.Item(3, 16).Value = 2.04
$ws.Cells.Item(3, 17).Value = 1.58
$ws.Cells.Item(3, 18).Value = 1.43
$ws.Cells.Item(3, 19).Value = 2.3
$ws.Cells.Item(3, 23).Value = 1.44
$ws.Cells.Item(3, 24).Value = 1000
$ws.Cells.Item(3, 25).Value = 1000
$ws.Cells.Item(3, 26).Value = 1000
$ws.Cells.Item(3, 28).Value = 1000
$ws.Cells.Item(3, 29).Value = 1000
$ws.Cells.Item(3, 30).Value = 1000
$ws.Cells.Item(3, 31).Value = 1000
$ws.Cells.Item(3, 32).Value = 1000
$ws.Cells.Item(3, 33).Value = 1000
$ws.Cells.Item(3, 34).Value = 1000
$ws.Cells.Item(3, 37).Value = 1000
$ws.Cells.Item(3, 40).Value = 1000
$ws.Cells.Item(3, 41).Value = 1000
$ws.Cells.Item(4, 6).Value = 1.56
$ws.Cells.Item(4, 7).Value = 1.68
$ws.Cells.Item(4, 8).Value = 4.8
$ws.Cells.Item(4, 9).Value = 8
$ws.Cells.Item(4, 11).Value = 4.7
$ws.Cells.Item(4, 14).Value = 2
$ws.Cells.Item(4, 15).Value = 1.28
$ws.Cells.Item(4, 17).Value = 1.28
$ws.Cells.Item(4, 18).Value = 1.34
$ws.Cells.Item(4, 19).Value = 2.68
$ws.Cells.Item(4, 22).Value = 1.14
$ws.Cells.Item(4, 23).Value = 2.28
$ws.Cells.Item(4, 25).Value = 1000
$ws.Cells.Item(4, 26).Value = 1000
$ws.Cells.Item(4, 28).Value = 1000
$ws.Cells.Item(4, 29).Value = 1000
$ws.Cells.Item(4, 30).Value = 1000
$ws.Cells.Item(4, 32).Value = 1000
$ws.Cells.Item(4, 33).Value = 1000
$ws.Cells.Item(4, 34).Value = 1000
$ws.Cells.Item(4, 36).Value = 1000
$ws.Cells.Item(4, 37).Value = 1000
$ws.Cells.Item(4, 38).Value = 1000
$ws.Cells.Item(5, 6).Value = 4.7
$ws.Cells.Item(5, 7).Value = 5.6
$ws.Cells.Item(5, 8).Value = 1.8
$ws.Cells.Item(5, 9).Value = 1.94
$ws.Cells.Item(5, 10).Value = 3.6
$ws.Cells.Item(5, 11).Value = 3.95
$ws.Cells.Item(5, 12).Value = 1.45
$ws.Cells.Item(5, 14).Value = 3.05
$ws.Cells.Item(5, 15).Value = 1.39
$ws.Cells.Item(5, 16).Value = 1.71
$ws.Cells.Item(5, 17).Value = 2.14
$ws.Cells.Item(5, 19).Value = 4
$ws.Cells.Item(5, 20).Value = 1.87
$ws.Cells.Item(5, 21).Value = 1.74
$ws.Cells.Item(5, 22).Value = 2.06
$ws.Cells.Item(5, 23).Value = 1.22
$ws.Cells.Item(5, 35).Value = 980
$ws.Cells.Item(6, 6).Value = 1.38
$ws.Cells.Item(6, 7).Value = 1.63
$ws.Cells.Item(6, 8).Value = 7.2
$ws.Cells.Item(6, 10).Value = 3.8
$ws.Cells.Item(6, 11).Value = 7.4
$ws.Cells.Item(6, 16).Value = 1.74
$ws.Cells.Item(6, 17).Value = 1.92
$ws.Cells.Item(7, 6).Value = 1.25
$ws.Cells.Item(7, 7).Value = 2.24
$ws.Cells.Item(7, 10).Value = 2.82
$ws.Cells.Item(7, 11).Value = 1000
$ws.Cells.Item(7, 16).Value = 1.25
$ws.Cells.Item(7, 17).Value = 1.01
$ws.Cells.Item(8, 6).Value = 1.09
$ws.Cells.Item(8, 9).Value = 1000
$ws.Cells.Item(8, 11).Value = 1000
$ws.Cells.Item(8, 16).Value = 1.15
$ws.Cells.Item(8, 17).Value = 1.01
$ws.Cells.Item(9, 6).Value = 2.18
$ws.Cells.Item(9, 8).Value = 2.98
$ws.Cells.Item(9, 10).Value = 2.44
$ws.Cells.Item(9, 16).Value = 1.84
$ws.Cells.Item(9, 17).Value = 1.72
$ws.Cells.Item(10, 16).Value = 1.62
$ws.Cells.Item(10, 18).Value = 1.23
$ws.Cells.Item(10, 20).Value = 2.02
$ws.Cells.Item(10, 24).Value = 9.4
$ws.Cells.Item(10, 25).Value = 10.5
$ws.Cells.Item(10, 26).Value = 22
$ws.Cells.Item(10, 27).Value = 80
$ws.Cells.Item(10, 28).Value = 8.800000000000001
$ws.Cells.Item(10, 30).Value = 15
$ws.Cells.Item(10, 31).Value = 48
$ws.Cells.Item(10, 32).Value = 16.5
$ws.Cells.Item(10, 33).Value = 13
$ws.Cells.Item(10, 34).Value = 20
$ws.Cells.Item(10, 35).Value = 70
$ws.Cells.Item(10, 36).Value = 44
$ws.Cells.Item(10, 37).Value = 36
$ws.Cells.Item(10, 38).Value = 70
$ws.Cells.Item(10, 39).Value = 170
$ws.Cells.Item(10, 40).Value = 38
$ws.Cells.Item(10, 41).Value = 60
$ws.Cells.Item(11, 7).Value = 980
$ws.Cells.Item(11, 10).Value = 3.15
$ws.Cells.Item(11, 16).Value = 1.6
$ws.Cells.Item(11, 17).Value = 2.02
$ws.Cells.Item(14, 6).Value = 3.65
$ws.Cells.Item(14, 8).Value = 1.73
$ws.Cells.Item(14, 16).Value = 1.99
$ws.Cells.Item(14, 17).Value = 1.6
$ws.Cells.Item(15, 6).Value = 2.16
$ws.Cells.Item(15, 7).Value = 2.36
$ws.Cells.Item(15, 8).Value = 3.7
$ws.Cells.Item(15, 9).Value = 4.3
$ws.Cells.Item(15, 10).Value = 3.1
$ws.Cells.Item(15, 11).Value = 3.7
$ws.Cells.Item(15, 17).Value = 2.36
$ws.Cells.Item(16, 6).Value = 1.87
$ws.Cells.Item(16, 7).Value = 2.08
$ws.Cells.Item(16, 8).Value = 3.4
$ws.Cells.Item(16, 9).Value = 4.3
$ws.Cells.Item(16, 10).Value = 4.3
$ws.Cells.Item(16, 11).Value = 5.8
$ws.Cells.Item(16, 16).Value = 3.15
$ws.Cells.Item(16, 17).Value = 1.37
$ws.Cells.Item(17, 6).Value = 3.45
$ws.Cells.Item(17, 7).Value = 5.4
$ws.Cells.Item(17, 8).Value = 1.75
$ws.Cells.Item(17, 9).Value = 2.2
$ws.Cells.Item(17, 10).Value = 2.8
$ws.Cells.Item(17, 11).Value = 7.8
$ws.Cells.Item(17, 16).Value = 2.2
$ws.Cells.Item(17, 17).Value = 1.54
$ws.Cells.Item(19, 6).Value = 2.26
$ws.Cells.Item(19, 7).Value = 2.94
$ws.Cells.Item(19, 8).Value = 2.96
$ws.Cells.Item(19, 9).Value = 4.2
$ws.Cells.Item(19, 10).Value = 3.05
$ws.Cells.Item(20, 6).Value = 2.6
$ws.Cells.Item(20, 7).Value = 2.98
$ws.Cells.Item(20, 8).Value = 2.76
$ws.Cells.Item(20, 10).Value = 3.05
$ws.Cells.Item(20, 11).Value = 3.6
$ws.Cells.Item(20, 16).Value = 1.73
$ws.Cells.Item(20, 17).Value = 1.97
$ws.Cells.Item(21, 14).Value = 2.96
$ws.Cells.Item(21, 16).Value = 1.62
$ws.Cells.Item(21, 20).Value = 2
$ws.Cells.Item(21, 24).Value = 9.199999999999999
$ws.Cells.Item(21, 25).Value = 9.6
$ws.Cells.Item(21, 26).Value = 19
$ws.Cells.Item(21, 27).Value = 65
$ws.Cells.Item(21, 28).Value = 9.4
$ws.Cells.Item(21, 30).Value = 14
$ws.Cells.Item(21, 31).Value = 42
$ws.Cells.Item(21, 33).Value = 13.5
$ws.Cells.Item(21, 34).Value = 19.5
$ws.Cells.Item(21, 35).Value = 60
$ws.Cells.Item(21, 36).Value = 50
$ws.Cells.Item(21, 37).Value = 40
$ws.Cells.Item(21, 38).Value = 75
$ws.Cells.Item(21, 39).Value = 170
$ws.Cells.Item(21, 40).Value = 44
$ws.Cells.Item(21, 41).Value = 46
$ws.Cells.Item(22, 7).Value = 1.95
$ws.Cells.Item(22, 16).Value = 2.76
$ws.Cells.Item(22, 17).Value = 1.54
$ws.Cells.Item(22, 26).Value = 36
$ws.Cells.Item(22, 27).Value = 80
$ws.Cells.Item(22, 28).Value = 15
$ws.Cells.Item(22, 31).Value = 42
$ws.Cells.Item(22, 35).Value = 36
$ws.Cells.Item(22, 38).Value = 25
$ws.Cells.Item(22, 39).Value = 55
$ws.Cells.Item(22, 40).Value = 8
$ws.Cells.Item(22, 41).Value = 29
$ws.Cells.Item(23, 6).Value = 2.32
$ws.Cells.Item(23, 15).Value = 1.42
$ws.Cells.Item(23, 16).Value = 1.75
$ws.Cells.Item(23, 17).Value = 2.3
$ws.Cells.Item(23, 19).Value = 4.3
$ws.Cells.Item(23, 20).Value = 1.98
$ws.Cells.Item(23, 21).Value = 1.98
$ws.Cells.Item(23, 24).Value = 11.5
$ws.Cells.Item(23, 27).Value = 70
$ws.Cells.Item(23, 28).Value = 8.6
$ws.Cells.Item(23, 35).Value = 60
$ws.Cells.Item(23, 36).Value = 34
$ws.Cells.Item(24, 6).Value = 2.86
$ws.Cells.Item(24, 7).Value = 3.1
$ws.Cells.Item(24, 8).Value = 2.68
$ws.Cells.Item(24, 9).Value = 2.8
$ws.Cells.Item(24, 12).Value = 1.4
$ws.Cells.Item(24, 14).Value = 3.75
$ws.Cells.Item(24, 16).Value = 1.92
$ws.Cells.Item(24, 17).Value = 1.94
$ws.Cells.Item(24, 18).Value = 1.36
$ws.Cells.Item(24, 19).Value = 3.35
$ws.Cells.Item(24, 20).Value = 1.72
$ws.Cells.Item(24, 22).Value = 1.55
$ws.Cells.Item(24, 23).Value = 1.48
$ws.Cells.Item(24, 24).Value = 15.5
$ws.Cells.Item(24, 25).Value = 11.5
$ws.Cells.Item(24, 26).Value = 18.5
$ws.Cells.Item(24, 27).Value = 980
$ws.Cells.Item(24, 28).Value = 12.5
$ws.Cells.Item(24, 29).Value = 7.8
$ws.Cells.Item(24, 30).Value = 12.5
$ws.Cells.Item(24, 31).Value = 30
$ws.Cells.Item(24, 32).Value = 21
$ws.Cells.Item(24, 33).Value = 13.5
$ws.Cells.Item(24, 34).Value = 16.5
$ws.Cells.Item(24, 35).Value = 42
$ws.Cells.Item(24, 36).Value = 60
$ws.Cells.Item(24, 37).Value = 40
$ws.Cells.Item(24, 38).Value = 980
$ws.Cells.Item(24, 39).Value = 100
$ws.Cells.Item(24, 40).Value = 29
$ws.Cells.Item(24, 41).Value = 24
$ws.Cells.Item(25, 7).Value = 2.76
$ws.Cells.Item(25, 8).Value = 3.1
$ws.Cells.Item(25, 9).Value = 4.8
$ws.Cells.Item(25, 10).Value = 2.78
$ws.Cells.Item(25, 11).Value = 4.5
$ws.Cells.Item(25, 14).Value = 1.56
$ws.Cells.Item(25, 16).Value = 1.56
$ws.Cells.Item(25, 17).Value = 2.1
$ws.Cells.Item(25, 18).Value = 1.19
$ws.Cells.Item(25, 19).Value = 3.9
$ws.Cells.Item(25, 22).Value = 1.26
$ws.Cells.Item(25, 23).Value = 1.57
$ws.Cells.Item(25, 24).Value = 1000
$ws.Cells.Item(25, 25).Value = 1000
$ws.Cells.Item(25, 26).Value = 1000
$ws.Cells.Item(25, 27).Value = 1000
$ws.Cells.Item(25, 28).Value = 1000
$ws.Cells.Item(25, 29).Value = 1000
$ws.Cells.Item(25, 30).Value = 1000
$ws.Cells.Item(25, 31).Value = 1000
$ws.Cells.Item(25, 32).Value = 1000
$ws.Cells.Item(25, 33).Value = 1000
$ws.Cells.Item(25, 34).Value = 1000
$ws.Cells.Item(25, 35).Value = 1000
$ws.Cells.Item(25, 36).Value = 1000
$ws.Cells.Item(25, 37).Value = 1000
$ws.Cells.Item(25, 38).Value = 1000
$ws.Cells.Item(25, 39).Value = 1000
$ws.Cells.Item(25, 40).Value = 1000
$ws.Cells.Item(25, 41).Value = 1000
